# Generate Report for Handback
# Populates the "Latest Target File" / "Latest Handback File" / "Latest
# Handback DateTime" columns for the zh-cn and de-de localization status
# sheets now that handback has completed, flips the Status text, adds the
# "back to GitHub" hyperlinks on the new target-file cells, and widens the
# columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

$targetFileDisplay = "2c38df77-837e-43d9-b4c1-8276e797efe4.md"
$targetFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/edf86a49926fc6dce97a773cf7af225530c26936/e2e/2c38df77-837e-43d9-b4c1-8276e797efe4.md"

# --- Status column: handback is in sync with en-US now ---------------------
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn: target file / handback file / handback datetime ---------------
$zh.Hyperlinks.Add($zh.Range("I2"), $targetFileUrl, "", "", $targetFileDisplay) | Out-Null
$zh.Range("J2").Value = "2c38df77-837e-43d9-b4c1-8276e797efe4.820350ece49fe2ca1b5fc060f04ed6e44c428dec.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-27 04:57:30"

# --- de-de: target file / handback file / handback datetime ---------------
$de.Hyperlinks.Add($de.Range("I2"), $targetFileUrl, "", "", $targetFileDisplay) | Out-Null
$de.Range("J2").Value = "2c38df77-837e-43d9-b4c1-8276e797efe4.820350ece49fe2ca1b5fc060f04ed6e44c428dec.de-de.xlf"
$de.Range("K2").Value = "2016-08-27 04:57:38"

# --- widen columns that now hold longer content -----------------------------
$overview.Columns.Item(5).ColumnWidth = 29.144371396019366
$overview.Columns.Item(6).ColumnWidth = 29.144371396019366

$zh.Columns.Item(3).ColumnWidth = 29.144371396019366
$zh.Columns.Item(9).ColumnWidth = 39.166666666666664
$zh.Columns.Item(10).ColumnWidth = 39.166666666666664

$de.Columns.Item(3).ColumnWidth = 29.144371396019366
$de.Columns.Item(9).ColumnWidth = 39.166666666666664
$de.Columns.Item(10).ColumnWidth = 39.166666666666664
